$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 76

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01-07-2021"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = 27333
$ws.Cells.Item($row, 3).Value = 10020
$ws.Cells.Item($row, 4).Value = 1410
$ws.Cells.Item($row, 5).Value = 5925
$ws.Cells.Item($row, 6).Value = 2684
$ws.Cells.Item($row, 7).Value = 17314
$ws.Cells.Item($row, 8).Value = 11857
$ws.Cells.Item($row, 9).Value = 5456
